$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '58.875.31'
$ws.Range('E2').Value = '  -2.44%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.732.51'
$ws.Range('E3').Value = '  -5.63%  '
$ws.Range('E4').Value = '  -0.07%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '506.35'
$ws.Range('E5').Value = '  -3.38%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '141.86'
$ws.Range('E6').Value = '  +1.03%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.997'
$ws.Range('E7').Value = '  -0.35%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.532'
$ws.Range('E8').Value = '  -2.90%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '2.745.87'
$ws.Range('E9').Value = '  -5.22%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '6.06'
$ws.Range('E10').Value = '  +3.12%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.105'
$ws.Range('E11').Value = '  -1.11%  '
$ws.Range('E12').Value = '  -1.62%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.127'
$ws.Range('E13').Value = '  +1.31%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '3.212.56'
$ws.Range('E14').Value = '  -5.60%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '58.962.92'
$ws.Range('E15').Value = '  -2.43%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '21.74'
$ws.Range('E16').Value = '  -3.26%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.0000136'
$ws.Range('E17').Value = '  -2.46%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '2.735.86'
$ws.Range('E18').Value = '  -5.86%  '
$ws.Range('E19').Value = '  -3.02%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '11.03'
$ws.Range('E20').Value = '  -3.74%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '344.42'
$ws.Range('E21').Value = '  -3.92%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.26'
$ws.Range('E22').Value = '  -3.38%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.999'
$ws.Range('E23').Value = '  -0.11%  '
$ws.Range('E24').Value = '  -0.48%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '63.37'
$ws.Range('E25').Value = '  +0.38%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.428'
$ws.Range('E26').Value = '  -3.81%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.172'
$ws.Range('E27').Value = '  -4.35%  '
$ws.Range('E28').Value = '  -0.40%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '7.53'
$ws.Range('E29').Value = '  -2.51%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.0₃0834'
$ws.Range('E30').Value = '  -0.45%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.999'
$ws.Range('B32').Value = 'EthereumClassic'
$ws.Range('C32').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '19.22'
$ws.Range('E32').Value = '  -0.59%  '
$ws.Range('B33').Value = 'PancakeSwap'
$ws.Range('C33').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.60'
$ws.Range('E33').Value = '  -3.36%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '151.80'
$ws.Range('E34').Value = '  +1.15%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '5.43'
$ws.Range('E35').Value = '  -1.07%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '4.19'
$ws.Range('E36').Value = '  -1.66%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.953'
$ws.Range('E37').Value = '  -2.84%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.13'
$ws.Range('E38').Value = '  -4.16%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '36.03'
$ws.Range('E39').Value = '  -4.60%  '
$ws.Range('E40').Value = '  -5.25%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '3.55'
$ws.Range('E41').Value = '  -1.89%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.190.76'
$ws.Range('E42').Value = '  -5.46%  '
$ws.Range('B43').Value = 'Hedera'
$ws.Range('C43').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.0557'
$ws.Range('E43').Value = '  -1.94%  '
$ws.Range('B44').Value = 'FirstDigitalUSD'
$ws.Range('C44').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.996'
$ws.Range('E44').Value = '  -0.28%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.605'
$ws.Range('E45').Value = '  -5.52%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '19.08'
$ws.Range('E46').Value = '  -6.87%  '
$ws.Range('E47').Value = '  +0.40%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '4.75'
$ws.Range('E48').Value = '  -4.32%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0226'
$ws.Range('E49').Value = '  -2.26%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0887'
$ws.Range('E50').Value = '  -3.89%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '18.12'
$ws.Range('E51').Value = '  +0.25%  '
